# Update the "testdata" worksheet:
#  - the secret-key customer block (rows 3-6) gets new names ("...4561" -> "...12877")
#    and new hyperlinked emails ("...4561@gmail.com" -> "...12876@gmail.com")
#  - the invalid-secret-key customer block (rows 10-13) gets new names
#    ("...456" -> "...12876")
#  - the deleteCustomer id list shrinks from two rows (17-18) to a single row (17)
#    with a new customer id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete second id row first so its old shared string is
# dropped before the replacement strings are appended.
$ws.Rows("18:18").Delete() | Out-Null

# Row 3-6 emails (column D, hyperlinked)
$ws.Range("D3").Value = "ag12876@gmail.com"
$ws.Range("D4").Value = "sk12876@gmail.com"
$ws.Range("D5").Value = "jsj12876@gmail.com"
$ws.Range("D6").Value = "kmrr12876@gmail.com"

# Row 10-13 names (column A)
$ws.Range("A10").Value = "Kumar Kishan12876"
$ws.Range("A11").Value = "Asif Khan12876"
$ws.Range("A12").Value = "Maharaj Saxena12876"
$ws.Range("A13").Value = "Kumar rawat12876"

# Row 3-6 names (column A)
$ws.Range("A3").Value = "Archana Gupta12877"
$ws.Range("A4").Value = "Suman Kumari12877"
$ws.Range("A5").Value = "Jyoti Saxena12877"
$ws.Range("A6").Value = "Kumar rawat rathode12877"

# New customer id replacing the old single id row
$ws.Range("A17").Value = "cus_Hb286lXLYZuvlx"

# Match the author's final selection
$ws.Range("A17").Select() | Out-Null
